$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B30: change from text "3" to numeric 3
$ws.Range("B30").Value = 3

# Add new row 31 with data
$ws.Range("A31").Value = "Ying Tang"
$ws.Range("B31").Value = "'4"
$ws.Range("B31").ClearFormats()
$ws.Range("C31").Value = "elaborate"
$ws.Range("D31").Value = "ACK"
$ws.Range("E31").Value = "WRI"
$ws.Range("F31").Value = "9b81a0cf-ae6f-4476-b619-1b75e1becf94"
$ws.Range("G31").Value = "B1ae1lZRb_annotated.xlsx"
$ws.Range("H31").Value = "We will elaborate on this aspect in the final version of the paper."
